$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update semester/index values in column A for rows 102 through 301 to 1
$ws.Range("A102:A301").Value = 1

# Update the active window scroll position / selection to match the new view
$win = $excel.ActiveWindow
$win.ScrollRow = 273
$win.ScrollColumn = 1
$ws.Range("E298").Select()
